# Katalog guncellendi - Per 27.11.2025 11:29:11,88
# Adds 3 new products (BOLD BROTHERS 6033 KASE GOMLEK - Lacivert / Gri / Yesil)
# to the bottom of the product catalog on Sheet1 (rows 115-117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$price = "475 TL"
$category = "Gömlek"
$desc = "S-M-L-XL  Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$stock = "Var"

# Row 115 - Lacivert
$ws.Range("A115").Value = "BOLD BROTHERS 6033 KAŞE GÖMLEK LACİVERT"
$ws.Range("B115").Value = $price
$ws.Range("C115").Value = $category
$ws.Range("D115").Value = "KAŞEGÖMLEKLACİ.jpg"
$ws.Range("E115").Value = $desc
$ws.Range("F115").Value = $stock

# Row 116 - Gri
$ws.Range("A116").Value = "BOLD BROTHERS 6033 KAŞE GÖMLEK GRİ"
$ws.Range("B116").Value = $price
$ws.Range("C116").Value = $category
$ws.Range("D116").Value = "KAŞEGÖMLEKGRİ.jpg"
$ws.Range("E116").Value = $desc
$ws.Range("F116").Value = $stock

# Row 117 - Yesil
$ws.Range("A117").Value = "BOLD BROTHERS 6033 KAŞE GÖMLEK YEŞİL"
$ws.Range("B117").Value = $price
$ws.Range("C117").Value = $category
$ws.Range("D117").Value = "KAŞEGÖMLEKYEŞİL.jpg"
$ws.Range("E117").Value = $desc
$ws.Range("F117").Value = $stock

$ws.Range("B126").Select()
